$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Implementação de frases marcantes, ou de músicas temas da abertura do anime",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Músicas temas da abertura do anime",
    2)
